$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.908.45'
$ws.Range('E2').Value = '  +2.44%  '

$ws.Range('D3').Value = '2.997.00'
$ws.Range('E3').Value = '  +1.78%  '

$ws.Range('E4').Value = '  +0.04%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '560.14'
$ws.Range('E5').Value = '  +1.28%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '136.91'
$ws.Range('E6').Value = '  +3.96%  '

$ws.Range('E7').Value = '  -0.16%  '

$ws.Range('E8').Value = '  +1.63%  '

$ws.Range('D9').Value = '2.990.58'
$ws.Range('E9').Value = '  +1.81%  '

$ws.Range('E10').Value = '  +4.11%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.19'
$ws.Range('E11').Value = '  +7.92%  '

$ws.Range('E12').Value = '  +2.26%  '

$ws.Range('E13').Value = '  +3.61%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '33.57'
$ws.Range('E14').Value = '  +2.60%  '

$ws.Range('E15').Value = '  +2.27%  '

$ws.Range('B16').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C16').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D16').Value = '3.494.63'
$ws.Range('E16').Value = '  +1.87%  '

$ws.Range('B17').Value = 'Polkadot'
$ws.Range('C17').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.31'
$ws.Range('E17').Value = '  +8.09%  '

$ws.Range('D18').Value = '2.996.30'
$ws.Range('E18').Value = '  +1.89%  '

$ws.Range('D19').Value = '58.953.51'
$ws.Range('E19').Value = '  +2.51%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '427.84'
$ws.Range('E20').Value = '  +3.12%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.74'
$ws.Range('E21').Value = '  +5.01%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.723'
$ws.Range('E22').Value = '  +6.10%  '

$ws.Range('E23').Value = '  +2.44%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.29'
$ws.Range('E24').Value = '  +2.56%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '80.51'
$ws.Range('E25').Value = '  +1.89%  '

$ws.Range('E26').Value = '  +0.04%  '

$ws.Range('E27').Value = '  +0.12%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.19'
$ws.Range('E28').Value = '  +11.17%  '

$ws.Range('E29').Value = '  +2.50%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.77'
$ws.Range('E30').Value = '  +3.67%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '25.76'
$ws.Range('E31').Value = '  +2.65%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.03'
$ws.Range('E32').Value = '  -0.05%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0991'
$ws.Range('E33').Value = '  -2.31%  '

$ws.Range('B34').Value = 'Mantle'
$ws.Range('C34').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.998'
$ws.Range('E34').Value = '  +6.67%  '

$ws.Range('B35').Value = 'Filecoin'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.95'
$ws.Range('E35').Value = '  +6.00%  '

$ws.Range('D36').Value = '0.0₃0755'
$ws.Range('E36').Value = '  +10.25%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.08'
$ws.Range('E37').Value = '  -0.77%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '48.72'
$ws.Range('E38').Value = '  +0.46%  '

$ws.Range('E39').Value = '  +2.98%  '

$ws.Range('E40').Value = '  +7.19%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '397.56'
$ws.Range('E41').Value = '  +5.28%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0351'
$ws.Range('E42').Value = '  +0.65%  '

$ws.Range('D43').Value = '2.751.97'
$ws.Range('E43').Value = '  +3.45%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.108'
$ws.Range('E44').Value = '  -0.55%  '

$ws.Range('E45').Value = '  +5.24%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '35.37'
$ws.Range('E46').Value = '  +25.54%  '

$ws.Range('E47').Value = '  -0.02%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '123.33'
$ws.Range('E48').Value = '  +0.87%  '

$ws.Range('E49').Value = '  +1.05%  '

$ws.Range('E50').Value = '  +0.88%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '23.34'
$ws.Range('E51').Value = '  +0.02%  '
